# Updates Leve profit-calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with refreshed market-price figures from the scheduled data-sync runner.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3800.3333
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 4180.4
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 4180.4
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -4992.4
$ws.Range("H91").Value = 3800.3333
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 4180.4
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 4180.4
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -6988.4
$ws.Range("H103").Value = 429.25
$ws.Range("I103").Value = 357.875
$ws.Range("K103").Value = 1073.625
$ws.Range("M103").Value = -487.625
$ws.Range("H138").Value = 4507.34
$ws.Range("I138").Value = 1983.8572
$ws.Range("J138").Value = 4891.3477
$ws.Range("K138").Value = 5951.571599999999
$ws.Range("L138").Value = 14674.0431
$ws.Range("M138").Value = -811.5715999999993
$ws.Range("N138").Value = -24954.0431

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2085.0278
$ws.Range("I122").Value = 1205.6552
$ws.Range("K122").Value = 3616.9656
$ws.Range("M122").Value = -1166.9656

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 33595
$ws.Range("J74").Value = 33595
$ws.Range("L74").Value = 33595
$ws.Range("N74").Value = -35467
$ws.Range("H77").Value = 33595
$ws.Range("J77").Value = 33595
$ws.Range("L77").Value = 100785
$ws.Range("N77").Value = -110145
$ws.Range("H105").Value = 3832.6128
$ws.Range("I105").Value = 2968.7727
$ws.Range("J105").Value = 5944.222
$ws.Range("K105").Value = 2968.7727
$ws.Range("L105").Value = 5944.222
$ws.Range("M105").Value = -1221.7727
$ws.Range("N105").Value = -9438.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 201.77777
$ws.Range("I7").Value = 64.5
$ws.Range("J7").Value = 311.6
$ws.Range("K7").Value = 64.5
$ws.Range("L7").Value = 311.6
$ws.Range("M7").Value = 48.5
$ws.Range("N7").Value = -537.6
$ws.Range("H16").Value = 35718140
$ws.Range("I16").Value = 50002588
$ws.Range("K16").Value = 50002588
$ws.Range("M16").Value = -50002301
$ws.Range("H31").Value = 4467.5415
$ws.Range("I31").Value = 1955.9
$ws.Range("J31").Value = 5128.5
$ws.Range("K31").Value = 1955.9
$ws.Range("L31").Value = 5128.5
$ws.Range("M31").Value = -1660.9
$ws.Range("N31").Value = -5718.5
$ws.Range("H34").Value = 4467.5415
$ws.Range("I34").Value = 1955.9
$ws.Range("J34").Value = 5128.5
$ws.Range("K34").Value = 1955.9
$ws.Range("L34").Value = 5128.5
$ws.Range("M34").Value = -1753.9
$ws.Range("N34").Value = -5532.5
$ws.Range("H113").Value = 35718140
$ws.Range("I113").Value = 50002588
$ws.Range("K113").Value = 50002588
$ws.Range("M113").Value = -50000418
$ws.Range("H132").Value = 5152.9375
$ws.Range("I132").Value = 4470.5854
$ws.Range("K132").Value = 13411.7562
$ws.Range("M132").Value = -10881.7562
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11000193
$ws.Range("I4").Value = 11000193
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 33000579
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -33000467
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 54.333332
$ws.Range("I2").Value = 47.125
$ws.Range("K2").Value = 47.125
$ws.Range("M2").Value = 65.875
$ws.Range("H80").Value = 66679596
$ws.Range("I80").Value = 125003430
$ws.Range("K80").Value = 125003430
$ws.Range("M80").Value = -125002432
$ws.Range("H83").Value = 66679596
$ws.Range("I83").Value = 125003430
$ws.Range("K83").Value = 625017150
$ws.Range("M83").Value = -625012158
$ws.Range("H102").Value = 942549.3
$ws.Range("I102").Value = 1670529.1
$ws.Range("J102").Value = 6575.3335
$ws.Range("K102").Value = 1670529.1
$ws.Range("L102").Value = 6575.3335
$ws.Range("M102").Value = -1668907.1
$ws.Range("N102").Value = -9819.3335
$ws.Range("H122").Value = 4543.5
$ws.Range("I122").Value = 3532.1667
$ws.Range("J122").Value = 5554.8335
$ws.Range("K122").Value = 10596.5001
$ws.Range("L122").Value = 16664.5005
$ws.Range("M122").Value = -8146.500100000001
$ws.Range("N122").Value = -21564.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3160.3333
$ws.Range("I16").Value = 918.3333
$ws.Range("K16").Value = 918.3333
$ws.Range("M16").Value = -748.3333
$ws.Range("H22").Value = 765.44446
$ws.Range("I22").Value = 482.6
$ws.Range("J22").Value = 1119
$ws.Range("K22").Value = 482.6
$ws.Range("L22").Value = 1119
$ws.Range("M22").Value = -187.6
$ws.Range("N22").Value = -1709
$ws.Range("H27").Value = 765.44446
$ws.Range("I27").Value = 482.6
$ws.Range("J27").Value = 1119
$ws.Range("K27").Value = 482.6
$ws.Range("L27").Value = 1119
$ws.Range("M27").Value = -375.6
$ws.Range("N27").Value = -1333
$ws.Range("H40").Value = 3613.3704
$ws.Range("I40").Value = 3204.1333
$ws.Range("K40").Value = 3204.1333
$ws.Range("M40").Value = -3068.1333
$ws.Range("H55").Value = 3966.475
$ws.Range("I55").Value = 3015.4482
$ws.Range("K55").Value = 3015.4482
$ws.Range("M55").Value = -2842.4482
$ws.Range("H101").Value = 31502.5
$ws.Range("J101").Value = 31502.5
$ws.Range("L101").Value = 31502.5
$ws.Range("N101").Value = -37992.5
$ws.Range("H122").Value = 7508.1665
$ws.Range("I122").Value = 9393.5
$ws.Range("K122").Value = 28180.5
$ws.Range("M122").Value = -25730.5
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 5281.9
$ws.Range("I132").Value = 4964.607
$ws.Range("K132").Value = 14893.821
$ws.Range("M132").Value = -12363.821

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 41000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 41000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 41000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -42020
$ws.Range("H62").Value = 10563
$ws.Range("I62").Value = 9051
$ws.Range("J62").Value = 12679.8
$ws.Range("K62").Value = 9051
$ws.Range("L62").Value = 12679.8
$ws.Range("M62").Value = -8427
$ws.Range("N62").Value = -13927.8
$ws.Range("H65").Value = 10563
$ws.Range("I65").Value = 9051
$ws.Range("J65").Value = 12679.8
$ws.Range("K65").Value = 45255
$ws.Range("L65").Value = 63399
$ws.Range("M65").Value = -42135
$ws.Range("N65").Value = -69639
$ws.Range("H100").Value = 2749.1875
$ws.Range("I100").Value = 1589.3334
$ws.Range("J100").Value = 4240.4287
$ws.Range("K100").Value = 3178.6668
$ws.Range("L100").Value = 8480.8574
$ws.Range("M100").Value = -2637.6668
$ws.Range("N100").Value = -9562.8574
$ws.Range("H122").Value = 2649.05
$ws.Range("I122").Value = 2511.9375
$ws.Range("J122").Value = 3197.5
$ws.Range("K122").Value = 7535.8125
$ws.Range("L122").Value = 9592.5
$ws.Range("M122").Value = -5085.8125
$ws.Range("N122").Value = -14492.5
$ws.Range("H132").Value = 6843.8237
$ws.Range("I132").Value = 4988
$ws.Range("J132").Value = 12875.25
$ws.Range("K132").Value = 14964
$ws.Range("L132").Value = 38625.75
$ws.Range("M132").Value = -12434
$ws.Range("N132").Value = -43685.75
